$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column (R) of data for year 2021, mirroring the existing
# column Q (year 2020) formatting.

# 1) Copy the formatting from column Q (rows 3-8) into column R so the
#    new cells inherit the same number formats / fonts / borders.
$ws.Range("Q3:Q8").Copy() | Out-Null
$ws.Range("R3:R8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# 2) Fill in the new column's values / formulas.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Formula = "=R6/R8*100000"
$ws.Range("R5").Formula = "=R7/R8*100000"
$ws.Range("R6").Value = 312
$ws.Range("R7").Value = 1910
$ws.Range("R8").Value = 4409166

# 3) Update the view: scroll back to the default top-left cell and move
#    the active selection to R15.
$ws.Range("R15").Select() | Out-Null
